# Account Payable.xlsx - add Journal Entry test-case sheets and extend the
# Invoice sheet with Close/Post transaction + Paytype + BankAccount columns;
# trim BankDetail back down to two rows and mirror its data onto a new
# "Sheet1" tab.

$wb = $excel.ActiveWorkbook

$wsInvoice = $wb.Worksheets.Item(1)
$wsBank    = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Add the three new worksheets at the end, in tab order.
# ---------------------------------------------------------------------
$wsSheet1          = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSheet1.Name     = "Sheet1"

$wsJournalEntriesD      = $wb.Worksheets.Add($null, $wsSheet1)
$wsJournalEntriesD.Name = "JournalEntriesD"

$wsJournalEntries      = $wb.Worksheets.Add($null, $wsJournalEntriesD)
$wsJournalEntries.Name = "JournalEntries"

# ---------------------------------------------------------------------
# 2. Invoice: new columns D:F (Close/Post transaction + Paytype)
# ---------------------------------------------------------------------
$wsInvoice.Range("D1").Value = "CloseTransaction"
$wsInvoice.Range("E1").Value = "PostTransaction"
$wsInvoice.Range("D2").Value = "Y"
$wsInvoice.Range("F1").Value = "Paytype"
$wsInvoice.Range("F2").Value = "EFT"
$wsInvoice.Range("E2").Value = "Y"
$wsInvoice.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$wsInvoice.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 3. JournalEntriesD / JournalEntries headers
# ---------------------------------------------------------------------
$wsJournalEntriesD.Range("A1").Value = "Account"
$wsJournalEntriesD.Range("B1").Value = "DebitAmount"
$wsJournalEntriesD.Range("C1").Value = "CreditAmount"

$wsJournalEntries.Range("A1").Value = "Reversal"

# JournalEntriesD data rows
$wsJournalEntriesD.Range("A2").Value = "6622 (Cash Account)"
$wsJournalEntriesD.Range("B2").Value = 100
$wsJournalEntriesD.Range("C2").Value = 0

$wsJournalEntriesD.Range("A3").Value = "1100 (Accounts Receivable Trade)"
$wsJournalEntriesD.Range("A3").Interior.Pattern = -4142   # xlPatternNone
$wsJournalEntriesD.Range("B3").Value = 0
$wsJournalEntriesD.Range("C3").Value = 100

$wsJournalEntriesD.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$wsJournalEntriesD.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$wsJournalEntriesD.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$wsJournalEntriesD.PageSetup.Orientation = 1

# JournalEntries remaining cells
$wsJournalEntries.Range("B1").Value = "EntryDesc"
$wsJournalEntries.Range("B2").Value = "NPAUTOTest"
$wsJournalEntries.Range("B3").Value = "NPAUTOTestRev"

$wsJournalEntries.Range("A2").Value = "N"
$wsJournalEntries.Range("A2").HorizontalAlignment = -4108   # xlCenter
$wsJournalEntries.Range("A3").Value = "Y"
$wsJournalEntries.Range("A3").HorizontalAlignment = -4108   # xlCenter
$wsJournalEntries.Rows.Item(2).RowHeight = 17.25
$wsJournalEntries.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 4. Invoice: BankAccount column (written last so its shared string
#    lands after the JournalEntries* strings, same as the source file).
# ---------------------------------------------------------------------
$wsInvoice.Range("G1").Value = "BankAccount"
$wsInvoice.Range("G2").Value = "a9d41000000012iAAA"

# ---------------------------------------------------------------------
# 5. BankDetail: drop row 3 (now duplicated on the new "Sheet1" tab)
# ---------------------------------------------------------------------
$wsBank.Rows.Item(3).Delete() | Out-Null

# ---------------------------------------------------------------------
# 6. "Sheet1": mirror the original 3-row BankDetail content
# ---------------------------------------------------------------------
$wsSheet1.Range("A1").Value = "Bank Account"
$wsSheet1.Range("A2").Value = "a9d41000000012iAAA"
$wsSheet1.Range("A3").Value = "a9d1K0000004DGVQA2"
$wsSheet1.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 7. Selections / active sheet (order matters - last one wins)
# ---------------------------------------------------------------------
$wsSheet1.Columns.Item(1).Select()
$wsJournalEntriesD.Range("H15").Select()
$wsJournalEntries.Range("H14").Select()
$wsBank.Range("A1:A2").Select()
$wsInvoice.Range("E11").Select()
